$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 114, shifting rows 114:182 down to 115:183.
$ws.Rows.Item(114).Insert()

# Populate the new row 114 with the updated record.
$ws.Cells.Item(114, 1).Value = 10
$ws.Cells.Item(114, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(114, 3).Value = "La Araucanía"
$ws.Cells.Item(114, 4).Value = 44438
$ws.Cells.Item(114, 5).Value = 9
$ws.Cells.Item(114, 6).Value = 100114013
$ws.Cells.Item(114, 7).Value = "Zanahoria"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 240
$ws.Cells.Item(114, 11).Value = 5000
$ws.Cells.Item(114, 12).Value = 6000
$ws.Cells.Item(114, 13).Value = 5583
$ws.Cells.Item(114, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(114, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(114, 16).Value = 223
$ws.Cells.Item(114, 17).Value = 25
$ws.Cells.Item(114, 18).Value = "Hortaliza"
